# PowerShell Excel COM-interop script
# Updates the cryptos list (Price and Volume(1h) columns) with refreshed
# market data, including the TrustWalletToken/Aptos row swap (rows 45-46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.019.74"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "1.902.66"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'0.7449"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'242.39"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.3074"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").Value = "'25.58"
$ws.Range("E9").Value = "  -6.69%  "
$ws.Range("D10").Value = "'0.06907"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "'0.08035"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'0.7556"
$ws.Range("E12").Value = "  -2.62%  "
$ws.Range("D13").Value = "1.902.11"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "'5.237"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "'91.27"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "'6.198"
$ws.Range("E16").Value = "  +3.85%  "
$ws.Range("D17").Value = "30.029.82"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "'14.05"
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").Value = "'0.000007785"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "'237.07"
$ws.Range("E20").Value = "  -5.29%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "2.155.02"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'7.101"
$ws.Range("E24").Value = "  +6.83%  "
$ws.Range("D25").Value = "'9.348"
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("D26").Value = "'167.89"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "'18.79"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "'0.1277"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Value = "'2.050"
$ws.Range("E29").Value = "  -5.30%  "
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "'4.307"
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "'4.050"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").Value = "'0.05287"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "'1.284"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").Value = "'0.7389"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").Value = "'2.727"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").Value = "'6.260"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("D41").Value = "'0.4465"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").Value = "'72.70"
$ws.Range("E42").Value = "  -5.22%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'7.752"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'0.8318"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "'101.53"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "'9.819"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").Value = "2.055.95"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "'36.60"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").Value = "'0.05988"
